$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new "Price" value is a plain decimal number: left alone, Excel
# would auto-coerce the assigned text into a numeric cell (the source data is
# always stored as text, e.g. "0.675"). Mark these cells as Text first so the
# value is kept as a literal string, then clear the formatting we just applied
# so the cell reverts to the workbook default style (no explicit format), just
# like every other text cell in this sheet.
$textCells = @("D5","D6","D7","D9","D11","D12","D16","D17","D20","D21","D22","D23","D24","D27","D28","D29","D30","D31","D32","D33","D34","D37","D38","D40","D41","D42","D43","D44","D45","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cryptos list values
$ws.Range('D2').Value = '43.967.43'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').Value = '2.353.46'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '0.675'
$ws.Range('E5').Value = '  -3.88%  '
$ws.Range('D6').Value = '240.78'
$ws.Range('E6').Value = '  -1.29%  '
$ws.Range('D7').Value = '73.52'
$ws.Range('E7').Value = '  -1.52%  '
$ws.Range('D9').Value = '0.605'
$ws.Range('E9').Value = '  +1.89%  '
$ws.Range('E10').Value = '  -2.20%  '
$ws.Range('D11').Value = '59.14'
$ws.Range('E11').Value = '  +2.36%  '
$ws.Range('D12').Value = '33.67'
$ws.Range('E12').Value = '  +5.73%  '
$ws.Range('E13').Value = '  -2.13%  '
$ws.Range('E14').Value = '  -0.55%  '
$ws.Range('D15').Value = '2.703.51'
$ws.Range('E15').Value = '  -0.79%  '
$ws.Range('D16').Value = '16.45'
$ws.Range('E16').Value = '  -3.42%  '
$ws.Range('D17').Value = '0.913'
$ws.Range('E17').Value = '  -0.86%  '
$ws.Range('D18').Value = '2.354.18'
$ws.Range('E18').Value = '  -0.69%  '
$ws.Range('D19').Value = '43.830.68'
$ws.Range('E19').Value = '  -1.37%  '
$ws.Range('D20').Value = '0.0000103'
$ws.Range('E20').Value = '  -1.07%  '
$ws.Range('D21').Value = '6.73'
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').Value = '77.79'
$ws.Range('E22').Value = '  -1.66%  '
$ws.Range('D23').Value = '257.44'
$ws.Range('E23').Value = '  -0.66%  '
$ws.Range('D24').Value = '1.93'
$ws.Range('E24').Value = '  +16.25%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = '2.51'
$ws.Range('E27').Value = '  -2.96%  '
$ws.Range('D28').Value = '10.67'
$ws.Range('E28').Value = '  -1.66%  '
$ws.Range('D29').Value = '2.28'
$ws.Range('E29').Value = '  -1.12%  '
$ws.Range('D30').Value = '22.70'
$ws.Range('E30').Value = '  -0.30%  '
$ws.Range('D31').Value = '177.33'
$ws.Range('E31').Value = '  +1.49%  '
$ws.Range('D32').Value = '0.129'
$ws.Range('E32').Value = '  -0.21%  '
$ws.Range('D33').Value = '0.137'
$ws.Range('E33').Value = '  -0.33%  '
$ws.Range('D34').Value = '0.0758'
$ws.Range('E34').Value = '  -0.41%  '
$ws.Range('E35').Value = '  -3.28%  '
$ws.Range('E36').Value = '  +1.82%  '
$ws.Range('D37').Value = '3.82'
$ws.Range('E37').Value = '  -2.67%  '
$ws.Range('D38').Value = '6.46'
$ws.Range('E38').Value = '  -1.89%  '
$ws.Range('E39').Value = '  -3.33%  '
$ws.Range('D40').Value = '0.0279'
$ws.Range('E40').Value = '  +0.50%  '
$ws.Range('D41').Value = '67.96'
$ws.Range('E41').Value = '  +27.56%  '
$ws.Range('D42').Value = '5.22'
$ws.Range('E42').Value = '  +17.11%  '
$ws.Range('D43').Value = '0.112'
$ws.Range('E43').Value = '  +9.38%  '
$ws.Range('D44').Value = '9.28'
$ws.Range('E44').Value = '  +1.83%  '
$ws.Range('D45').Value = '19.13'
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('E46').Value = '  +1.52%  '
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('E48').Value = '  -0.53%  '
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = '99.66'
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').Value = '1.16'
$ws.Range('E51').Value = '  -3.82%  '

# Drop the temporary Text format so the cell style matches the rest of the sheet
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
